## Applies the commit "update formatting; insert solution into doc;"
##
## 1. Replaces the trailing empty paragraph at the end of the body with two
##    new centered, bold paragraphs containing the worked-problem answers:
##        k = 250000000000143
##        p = 3000000000001727
## 2. Updates the "Normal" paragraph style's default font color from
##    "auto" to RGB 00000A.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert the two solution lines into the last (empty) paragraph.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$kPara = $d.Paragraphs($n)
$kPara.Range.Text = "k = 250000000000143"
$kPara.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$pPara = $d.Paragraphs($n)
$pPara.Range.Text = "p = 3000000000001727"

# Center + bold both new paragraphs (normal and complex-script bold).
$kPara.Alignment = 1
$kPara.Range.Font.Bold = 1
$kPara.Range.Font.BoldBi = 1

$pPara.Alignment = 1
$pPara.Range.Font.Bold = 1
$pPara.Range.Font.BoldBi = 1

# ---------------------------------------------------------------------
# 2) Change the "Normal" style's font color from automatic to 00000A.
#    Word's Font.Color is a BGR-packed integer (0x00BBGGRR); 00000A ->
#    R=0x00 G=0x00 B=0x0A -> 0x0A0000 == 655360.
# ---------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Color = 655360
